$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 19 (season 2004) previously stored the em dash "-" as text in C19:E19,
# with C19 left in the default (unstyled) format while D19/E19 used the
# bordered numeric style. Replace the values with numeric zeros and give
# C19 the same formatting as D19/E19 so all three match.
$ws.Range("D19").Copy()
$ws.Range("C19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 0

# Reflect the active selection recorded in the sheet view.
$ws.Range("C19").Select()
